$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dadosDeAcesso")

# Add new row 6 with data matching the existing pattern (ID, nomeDeUsuario, email, senha)
$ws.Range("A6").Value = "ID_0008"
$ws.Range("B6").Value = "André Automatizador"
$ws.Range("C6").Value = "sem email"
$ws.Range("D6").Value = "automacaoteste"

# Update the selection to reflect the new active cell
$ws.Activate()
$ws.Range("D6").Select()
